$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = 0.230598549219447
$ws.Range("C2").Value = 0.7445381377257259
$ws.Range("D2").Value = 1.018876335886205
$ws.Range("E2").Value = 1.009394043912587
$ws.Range("F2").Value = 1.011190667123304
$ws.Range("G2").Value = 18

# Row 3 (Q1)
$ws.Range("B3").Value = 0.3045708550604296
$ws.Range("C3").Value = 0.3907590825865097
$ws.Range("D3").Value = 0.2095424611763801
$ws.Range("E3").Value = 0.4577580814976182
$ws.Range("F3").Value = 0.3522467123879903
$ws.Range("G3").Value = 17

# Row 4 (Q2)
$ws.Range("B4").Value = 0.2502419823634077
$ws.Range("C4").Value = 0.366077697295096
$ws.Range("D4").Value = 0.1747842588224781
$ws.Range("E4").Value = 0.4180720737175327
$ws.Range("F4").Value = 0.345891249322188
$ws.Range("G4").Value = 16

# Row 5 (Q3)
$ws.Range("B5").Value = 0.3395964577034516
$ws.Range("C5").Value = 0.406420661795407
$ws.Range("D5").Value = 0.2182909649593746
$ws.Range("E5").Value = 0.4672161865340012
$ws.Range("F5").Value = 0.3321443493336289
$ws.Range("G5").Value = 15

# Row 6 (Q4)
$ws.Range("B6").Value = 0.3836894047943985
$ws.Range("C6").Value = 0.3861383846528136
$ws.Range("D6").Value = 0.2113509978471021
$ws.Range("E6").Value = 0.4597292658153297
$ws.Range("F6").Value = 0.2628055933924591
$ws.Range("G6").Value = 14

# Row 7 (Q5)
$ws.Range("B7").Value = 0.3925847413104661
$ws.Range("C7").Value = 0.4061565136574045
$ws.Range("D7").Value = 0.2245340717459249
$ws.Range("E7").Value = 0.4738502630007974
$ws.Range("F7").Value = 0.2761863507775066
$ws.Range("G7").Value = 13

# Row 8 (Q6)
$ws.Range("B8").Value = 0.4359548671747835
$ws.Range("C8").Value = 0.4359548671747835
$ws.Range("D8").Value = 0.2345424034301451
$ws.Range("E8").Value = 0.4842957809336615
$ws.Range("F8").Value = 0.2202950679514644
$ws.Range("G8").Value = 12

# Row 9 (Q7)
$ws.Range("B9").Value = 0.4182767482978756
$ws.Range("C9").Value = 0.4182767482978756
$ws.Range("D9").Value = 0.2153166270150599
$ws.Range("E9").Value = 0.4640222268545548
$ws.Range("F9").Value = 0.2107066864939436
$ws.Range("G9").Value = 11

# Row 10 (Q8)
$ws.Range("B10").Value = 0.3824308908543818
$ws.Range("C10").Value = 0.3824308908543818
$ws.Range("D10").Value = 0.1807831973670264
$ws.Range("E10").Value = 0.4251860738159546
$ws.Range("F10").Value = 0.1958735734184745
$ws.Range("G10").Value = 10

# Row 11 (Q9)
$ws.Range("B11").Value = 0.3843625910566515
$ws.Range("C11").Value = 0.3843625910566515
$ws.Range("D11").Value = 0.1817866140615611
$ws.Range("E11").Value = 0.4263644146285676
$ws.Range("F11").Value = 0.195725609566047
$ws.Range("G11").Value = 9
